# Update the "IMPLEMENTATION" row (row 8) with the missing amount of hours
# for FERMI (C), GRILLO (D) and JERRY (E), and move the active selection
# to E8 to match where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 165
$ws.Range("D8").Value = 205
$ws.Range("E8").Value = 70

$ws.Range("E8").Select()
